$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "2024-09-25T18:04:30Z"
$ws.Range("B6").Value = "temperature"
$ws.Range("C6").Value = "'28"
$ws.Range("D6").Value = "N/A"
$ws.Range("E6").Value = "N/A"
$ws.Range("F6").Value = "N/A"
